$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 111: finish the "Bitwise ORs of Subarrays" (#898) entry that
# previously only had columns A and B filled in. ---
$ws.Range("C111").Value = "#array #bit-minipulation #dynamic-programming"
$ws.Range("D111").Value = "medium"
$ws.Range("E111").Value = 0
$ws.Range("F111").Value = 1
$ws.Range("G111").Value = 30
$ws.Range("H111").Value = 45869
$ws.Range("I111").Value = 45869

# --- Row 112: LeetCode 118 - Pascal's Triangle ---
$ws.Range("A112").Value = 118
$ws.Range("B112").Value = "Pascal's Triangle"
$ws.Range("C112").Value = "#array #dynamic-programming "
$ws.Range("D112").Value = "easy"
$ws.Range("E112").Value = 1
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 6
$ws.Range("H112").Value = 45870
$ws.Range("I112").Value = 45870

# --- Row 113: LeetCode 314 - Binary Tree Vertical Order Traversal ---
$ws.Range("A113").Value = 314
$ws.Range("B113").Value = "Binary Tree Vertical Order Traversal"
$ws.Range("C113").Value = "#tree #binary-tree #hash-table #bfs "
$ws.Range("D113").Value = "medium"
$ws.Range("E113").Value = 1
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 23
$ws.Range("H113").Value = 45870
$ws.Range("I113").Value = 45870

# Match the formatting used by the other data rows (centered alignment,
# wrapped tag/name text, date number format) by copying the formats from
# the row above, column group by column group.
$ws.Range("A110").Copy()
$ws.Range("A111:A113").PasteSpecial(-4122)
$ws.Range("B110:C110").Copy()
$ws.Range("B111:C113").PasteSpecial(-4122)
$ws.Range("D110:G110").Copy()
$ws.Range("D111:G113").PasteSpecial(-4122)
$ws.Range("H110:I110").Copy()
$ws.Range("H111:I113").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights (auto-grown because of the wrapped Tags/Name text).
$ws.Rows.Item(111).RowHeight = 68
$ws.Rows.Item(112).RowHeight = 34
$ws.Rows.Item(113).RowHeight = 34

$ws.Range("D115").Select()
